$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples_retained")

# Insert a new column before column B ("Discourse or standalone"),
# shifting the existing B:J data right to C:K.
$ws.Columns("B:B").Insert()

# Fill in the previously-empty "tess" row (row 26) with its dataset details
# (new shared-string values are added in the same order the original author
# entered them, so the resulting shared-strings table lines up).
$ws.Range("C26").Value = "acted"
$ws.Range("D26").Value = 800
$ws.Range("E26").Value = 1600
$ws.Range("F26").Value = 400
$ws.Range("G26").Value = "English"
$ws.Range("I26").Value = 2
$ws.Range("K26").Value = "Toronto English"
$ws.Range("J26").Value = "anger, disgust, fear, happiness, pleasant surprise, sadness, neutral"

# New column header + width (matches col width of column A: 19, customWidth, no bestFit).
$ws.Range("B1").Value = "Discourse or standalone"
$ws.Columns("B:B").ColumnWidth = 18.17

# Clear the leftover cell-selection marker left over from editing in Excel.
$ws.Range("A1").Select()
